$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title cell (A1) ---
$ws.Range("A1").Value = "ElHadar-PLC  Statement of cash flow  From Fri Jan 01 1999 To Sat Jan 01 2022"

# --- Header row (row 4): Account -> (blank), (blank) -> Account Code, Total stays ---
$ws.Range("A4").Formula = '=""'
$ws.Range("B4").Value = "Account Code"
$ws.Range("C4").Value = "Total"

# --- Row 5: blank spacer row ---
$ws.Range("A5").Formula = '=""'

# --- Row 6 ---
$ws.Range("A6").Value = "Beginning Cash Balance"

# --- Row 8 ---
$ws.Range("A8").Value = "IDK"

# --- Row 10 ---
$ws.Range("A10").Value = " HELLO(HELLO)"
$ws.Range("B10").Formula = '=""'
$ws.Range("C10").Formula = '="-988"'

# --- Row 12 (totals row) ---
$ws.Range("A12").Formula = '=""'
$ws.Range("B12").Formula = '=""'
$ws.Range("C12").Formula = '="-988"'

# --- Merge region shrinks from A1:D3 to A1:C3 ---
$ws.Range("A1:D3").UnMerge()
$ws.Range("A1:C3").Merge()
